$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 6758259
$ws.Cells.Item(137, 9).Value = 1265.22
$ws.Cells.Item(137, 10).Value = 20835328
$ws.Cells.Item(137, 11).Value = 3795.66
$ws.Cells.Item(137, 12).Value = 62505984
$ws.Cells.Item(137, 13).Value = -1245.66
$ws.Cells.Item(137, 14).Value = -62511084

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 104694.4
$ws.Cells.Item(2, 9).Value = 146134.86
$ws.Cells.Item(2, 10).Value = 8000
$ws.Cells.Item(2, 11).Value = 146134.86
$ws.Cells.Item(2, 12).Value = 8000
$ws.Cells.Item(2, 13).Value = -146021.86
$ws.Cells.Item(2, 14).Value = -8226

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 176.66667
$ws.Cells.Item(5, 9).Value = 176.66667
$ws.Cells.Item(5, 11).Value = 176.66667
$ws.Cells.Item(5, 13).Value = -64.66667000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1597.1724
$ws.Cells.Item(61, 9).Value = 1374.9524
$ws.Cells.Item(61, 10).Value = 2180.5
$ws.Cells.Item(61, 11).Value = 1374.9524
$ws.Cells.Item(61, 12).Value = 2180.5
$ws.Cells.Item(61, 13).Value = -1162.9524
$ws.Cells.Item(61, 14).Value = -2604.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 1085.641
$ws.Cells.Item(74, 9).Value = 1352.3889
$ws.Cells.Item(74, 10).Value = 857
$ws.Cells.Item(74, 11).Value = 1352.3889
$ws.Cells.Item(74, 12).Value = 857
$ws.Cells.Item(74, 13).Value = -478.3888999999999
$ws.Cells.Item(74, 14).Value = -2605

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 1085.641
$ws.Cells.Item(77, 9).Value = 1352.3889
$ws.Cells.Item(77, 10).Value = 857
$ws.Cells.Item(77, 11).Value = 6761.9445
$ws.Cells.Item(77, 12).Value = 4285
$ws.Cells.Item(77, 13).Value = -2393.9445
$ws.Cells.Item(77, 14).Value = -13021

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 104694.4
$ws.Cells.Item(116, 9).Value = 146134.86
$ws.Cells.Item(116, 10).Value = 8000
$ws.Cells.Item(116, 11).Value = 146134.86
$ws.Cells.Item(116, 12).Value = 8000
$ws.Cells.Item(116, 13).Value = -143840.86
$ws.Cells.Item(116, 14).Value = -12588

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 894.3043
$ws.Cells.Item(122, 9).Value = 903.7895
$ws.Cells.Item(122, 11).Value = 2711.3685
$ws.Cells.Item(122, 13).Value = -261.3685

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 2211.611
$ws.Cells.Item(132, 9).Value = 1740.3478
$ws.Cells.Item(132, 11).Value = 5221.0434
$ws.Cells.Item(132, 13).Value = -2691.0434

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 1597.1724
$ws.Cells.Item(136, 9).Value = 1374.9524
$ws.Cells.Item(136, 10).Value = 2180.5
$ws.Cells.Item(136, 11).Value = 4124.857199999999
$ws.Cells.Item(136, 12).Value = 6541.5
$ws.Cells.Item(136, 13).Value = -1574.857199999999
$ws.Cells.Item(136, 14).Value = -11641.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 104694.4
$ws.Cells.Item(3, 9).Value = 146134.86
$ws.Cells.Item(3, 10).Value = 8000
$ws.Cells.Item(3, 11).Value = 146134.86
$ws.Cells.Item(3, 12).Value = 8000
$ws.Cells.Item(3, 13).Value = -146020.86
$ws.Cells.Item(3, 14).Value = -8228

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 176.66667
$ws.Cells.Item(4, 9).Value = 176.66667
$ws.Cells.Item(4, 11).Value = 176.66667
$ws.Cells.Item(4, 13).Value = -61.66667000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 2479.0833
$ws.Cells.Item(94, 9).Value = 2249.9092
$ws.Cells.Item(94, 11).Value = 2249.9092
$ws.Cells.Item(94, 13).Value = -1798.9092

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 5234.222
$ws.Cells.Item(58, 9).Value = 1058.8334
$ws.Cells.Item(58, 10).Value = 8574.533
$ws.Cells.Item(58, 11).Value = 1058.8334
$ws.Cells.Item(58, 12).Value = 8574.533
$ws.Cells.Item(58, 13).Value = -855.8334
$ws.Cells.Item(58, 14).Value = -8980.533

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 5234.222
$ws.Cells.Item(136, 9).Value = 1058.8334
$ws.Cells.Item(136, 10).Value = 8574.533
$ws.Cells.Item(136, 11).Value = 3176.5002
$ws.Cells.Item(136, 12).Value = 25723.599
$ws.Cells.Item(136, 13).Value = -626.5002
$ws.Cells.Item(136, 14).Value = -30823.599

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 1266.36
$ws.Cells.Item(68, 9).Value = 660.5405
$ws.Cells.Item(68, 10).Value = 1622.1587
$ws.Cells.Item(68, 11).Value = 1981.6215
$ws.Cells.Item(68, 12).Value = 4866.4761
$ws.Cells.Item(68, 13).Value = -1170.6215
$ws.Cells.Item(68, 14).Value = -6488.4761

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(71, 8).Value = 1266.36
$ws.Cells.Item(71, 9).Value = 660.5405
$ws.Cells.Item(71, 10).Value = 1622.1587
$ws.Cells.Item(71, 11).Value = 5944.8645
$ws.Cells.Item(71, 12).Value = 14599.4283
$ws.Cells.Item(71, 13).Value = -1888.8645
$ws.Cells.Item(71, 14).Value = -22711.4283

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 147925.98
$ws.Cells.Item(107, 9).Value = 313.8095
$ws.Cells.Item(107, 10).Value = 244796.47
$ws.Cells.Item(107, 11).Value = 941.4285
$ws.Cells.Item(107, 12).Value = 734389.41
$ws.Cells.Item(107, 13).Value = 978.5715
$ws.Cells.Item(107, 14).Value = -738229.41

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(32, 8).Value = 0
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 12).Value = 0
$ws.Cells.Item(32, 14).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 3505.5
$ws.Cells.Item(80, 9).Value = 3653.8235
$ws.Cells.Item(80, 10).Value = 3001.2
$ws.Cells.Item(80, 11).Value = 3653.8235
$ws.Cells.Item(80, 12).Value = 3001.2
$ws.Cells.Item(80, 13).Value = -2655.8235
$ws.Cells.Item(80, 14).Value = -4997.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(83, 8).Value = 3505.5
$ws.Cells.Item(83, 9).Value = 3653.8235
$ws.Cells.Item(83, 10).Value = 3001.2
$ws.Cells.Item(83, 11).Value = 18269.1175
$ws.Cells.Item(83, 12).Value = 15006
$ws.Cells.Item(83, 13).Value = -13277.1175
$ws.Cells.Item(83, 14).Value = -24990

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 870.6
$ws.Cells.Item(102, 9).Value = 732.2105
$ws.Cells.Item(102, 11).Value = 732.2105
$ws.Cells.Item(102, 13).Value = 889.7895

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 897.53845
$ws.Cells.Item(61, 9).Value = 888
$ws.Cells.Item(61, 10).Value = 950
$ws.Cells.Item(61, 11).Value = 888
$ws.Cells.Item(61, 12).Value = 950
$ws.Cells.Item(61, 13).Value = -686
$ws.Cells.Item(61, 14).Value = -1354

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 1691.9166
$ws.Cells.Item(82, 9).Value = 1730
$ws.Cells.Item(82, 10).Value = 1501.5
$ws.Cells.Item(82, 11).Value = 1730
$ws.Cells.Item(82, 12).Value = 1501.5
$ws.Cells.Item(82, 13).Value = -1369
$ws.Cells.Item(82, 14).Value = -2223.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(85, 8).Value = 1691.9166
$ws.Cells.Item(85, 9).Value = 1730
$ws.Cells.Item(85, 10).Value = 1501.5
$ws.Cells.Item(85, 11).Value = 1730
$ws.Cells.Item(85, 12).Value = 1501.5
$ws.Cells.Item(85, 13).Value = -482
$ws.Cells.Item(85, 14).Value = -3997.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 1680.2667
$ws.Cells.Item(93, 9).Value = 945.6667
$ws.Cells.Item(93, 10).Value = 2782.1667
$ws.Cells.Item(93, 11).Value = 945.6667
$ws.Cells.Item(93, 12).Value = 2782.1667
$ws.Cells.Item(93, 13).Value = 302.3333
$ws.Cells.Item(93, 14).Value = -5278.1667

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 897.53845
$ws.Cells.Item(113, 9).Value = 888
$ws.Cells.Item(113, 10).Value = 950
$ws.Cells.Item(113, 11).Value = 888
$ws.Cells.Item(113, 12).Value = 950
$ws.Cells.Item(113, 13).Value = 1282
$ws.Cells.Item(113, 14).Value = -5290

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 1872.3594
$ws.Cells.Item(136, 9).Value = 1129.425
$ws.Cells.Item(136, 10).Value = 3110.5833
$ws.Cells.Item(136, 11).Value = 3388.275
$ws.Cells.Item(136, 12).Value = 9331.749899999999
$ws.Cells.Item(136, 13).Value = -838.2749999999996
$ws.Cells.Item(136, 14).Value = -14431.7499

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 3000.5833
$ws.Cells.Item(96, 9).Value = 3086.1428
$ws.Cells.Item(96, 10).Value = 2880.8
$ws.Cells.Item(96, 11).Value = 3086.1428
$ws.Cells.Item(96, 12).Value = 2880.8
$ws.Cells.Item(96, 13).Value = -1713.1428
$ws.Cells.Item(96, 14).Value = -5626.8
